$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.737.68'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '2.552.77'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'302.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.67%  '
$ws.Range("D6").Value = "'98.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.25%  '
$ws.Range("D7").Value = "'0.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = "'0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").Value = "'36.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").Value = "'0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("E12").Value = '  +7.87%  '
$ws.Range("D13").Value = "'7.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.58%  '
$ws.Range("D14").Value = '2.531.62'
$ws.Range("E14").Value = '  -0.70%  '
$ws.Range("D15").Value = "'0.873"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = "'14.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.09%  '
$ws.Range("D17").Value = '42.786.55'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").Value = "'13.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.68%  '
$ws.Range("D19").Value = '0.0₃0983'
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").Value = "'71.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("D22").Value = "'253.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.63%  '
$ws.Range("D23").Value = "'2.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.33%  '
$ws.Range("D24").Value = "'2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.12%  '
$ws.Range("D25").Value = "'27.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.06%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = "'10.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = "'37.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.30%  '
$ws.Range("D29").Value = "'2.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("D30").Value = "'5.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").Value = "'155.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.65%  '
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("D34").Value = "'0.0804"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.32%  '
$ws.Range("D35").Value = "'3.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.66%  '
$ws.Range("D36").Value = "'26.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.08%  '
$ws.Range("D37").Value = "'18.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.98%  '
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("D39").Value = "'0.119"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("D42").Value = "'3.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.39%  '
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.067.72'
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = "'0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("D46").Value = "'87.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.94%  '
$ws.Range("D47").Value = "'9.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.96%  '
$ws.Range("D48").Value = '2.799.04'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").Value = "'74.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.09%  '
$ws.Range("D50").Value = "'103.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.02%  '
$ws.Range("E51").Value = '  +1.19%  '
